# Refresh odds/match data in the "Peru Liga 1" sheet (league database update).
# Several already-played fixtures had their scraped id/teams/odds corrected
# (mostly row-pair/row-triple re-matches), and several upcoming fixtures had
# their closing odds (oddH/oddD/oddA .. oddAHUnder, columns N:V) refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 210
$ws.Range("B210").Value = 6905578
$ws.Range("F210").Value = "AD Tarma"
$ws.Range("G210").Value = "Atletico Grau"
$ws.Range("I210").Value = 0
$ws.Range("J210").Value = "H"
$ws.Range("K210").Value = 1.75
$ws.Range("L210").Value = 3.6
$ws.Range("M210").Value = 4
$ws.Range("N210").Value = 1.571
$ws.Range("O210").Value = 4.2
$ws.Range("P210").Value = 5.75
$ws.Range("Q210").Value = -1
$ws.Range("R210").Value = 1.975
$ws.Range("S210").Value = 1.825
$ws.Range("U210").Value = 1.8
$ws.Range("V210").Value = 2
$ws.Range("W210").Value = 0.571
$ws.Range("X210").Value = -1
$ws.Range("Z210").Value = 0
$ws.Range("AA210").Value = 0
$ws.Range("AC210").Value = 1

# Row 211
$ws.Range("B211").Value = 6905571
$ws.Range("F211").Value = "FBC Melgar"
$ws.Range("G211").Value = "Sporting Cristal"
$ws.Range("I211").Value = 1
$ws.Range("J211").Value = "D"
$ws.Range("K211").Value = 2.1
$ws.Range("L211").Value = 3.4
$ws.Range("M211").Value = 3
$ws.Range("N211").Value = 1.75
$ws.Range("O211").Value = 3.8
$ws.Range("P211").Value = 4.75
$ws.Range("Q211").Value = -0.75
$ws.Range("R211").Value = 1.95
$ws.Range("S211").Value = 1.85
$ws.Range("U211").Value = 1.95
$ws.Range("V211").Value = 1.85
$ws.Range("W211").Value = -1
$ws.Range("X211").Value = 2.8
$ws.Range("Z211").Value = -1
$ws.Range("AA211").Value = 0.8500000000000001
$ws.Range("AC211").Value = 0.8500000000000001

# Row 305
$ws.Range("B305").Value = 7211640
$ws.Range("F305").Value = "UTC Cajamarca"
$ws.Range("G305").Value = "Sport Boys"
$ws.Range("H305").Value = 1
$ws.Range("I305").Value = 1
$ws.Range("J305").Value = "D"
$ws.Range("K305").Value = 1.615
$ws.Range("L305").Value = 3.75
$ws.Range("M305").Value = 5
$ws.Range("N305").Value = 1.5
$ws.Range("O305").Value = 4.2
$ws.Range("P305").Value = 6.5
$ws.Range("Q305").Value = -1
$ws.Range("R305").Value = 1.8
$ws.Range("S305").Value = 2.05
$ws.Range("T305").Value = 2.5
$ws.Range("U305").Value = 1.875
$ws.Range("V305").Value = 1.975
$ws.Range("W305").Value = -1
$ws.Range("X305").Value = 3.2
$ws.Range("Z305").Value = -1
$ws.Range("AA305").Value = 1.05
$ws.Range("AC305").Value = 0.9750000000000001

# Row 306
$ws.Range("B306").Value = 7211641
$ws.Range("F306").Value = "Sport Huancayo"
$ws.Range("G306").Value = "Deportivo Municipal"
$ws.Range("H306").Value = 2
$ws.Range("I306").Value = 0
$ws.Range("J306").Value = "H"
$ws.Range("K306").Value = 1.125
$ws.Range("L306").Value = 7
$ws.Range("M306").Value = 17
$ws.Range("N306").Value = 1.166
$ws.Range("O306").Value = 6.5
$ws.Range("P306").Value = 12
$ws.Range("Q306").Value = -2
$ws.Range("R306").Value = 1.775
$ws.Range("S306").Value = 2.025
$ws.Range("T306").Value = 3.5
$ws.Range("U306").Value = 1.9
$ws.Range("V306").Value = 1.9
$ws.Range("W306").Value = 0.1659999999999999
$ws.Range("X306").Value = -1
$ws.Range("Z306").Value = 0
$ws.Range("AA306").Value = 0
$ws.Range("AC306").Value = 0.8999999999999999

# Row 324
$ws.Range("B324").Value = 7302200
$ws.Range("F324").Value = "Carlos Manucci"
$ws.Range("G324").Value = "Deportivo Binacional"
$ws.Range("H324").Value = 3
$ws.Range("J324").Value = "H"
$ws.Range("K324").Value = 2
$ws.Range("L324").Value = 3.2
$ws.Range("M324").Value = 3.75
$ws.Range("O324").Value = 3.4
$ws.Range("P324").Value = 4.333
$ws.Range("R324").Value = 1.85
$ws.Range("S324").Value = 1.95
$ws.Range("T324").Value = 2.5
$ws.Range("U324").Value = 1.85
$ws.Range("V324").Value = 1.95
$ws.Range("W324").Value = 0.75
$ws.Range("Y324").Value = -1
$ws.Range("Z324").Value = 0.8500000000000001
$ws.Range("AA324").Value = -1
$ws.Range("AB324").Value = 0.8500000000000001
$ws.Range("AC324").Value = -1

# Row 325
$ws.Range("B325").Value = 7302796
$ws.Range("F325").Value = "Sport Huancayo"
$ws.Range("G325").Value = "Sport Boys"
$ws.Range("H325").Value = 1
$ws.Range("I325").Value = 0
$ws.Range("K325").Value = 1.727
$ws.Range("L325").Value = 3.75
$ws.Range("M325").Value = 4.333
$ws.Range("N325").Value = 1.25
$ws.Range("O325").Value = 5.25
$ws.Range("P325").Value = 10
$ws.Range("Q325").Value = -1.75
$ws.Range("R325").Value = 1.925
$ws.Range("S325").Value = 1.875
$ws.Range("T325").Value = 3
$ws.Range("U325").Value = 1.875
$ws.Range("V325").Value = 1.925
$ws.Range("W325").Value = 0.25
$ws.Range("Z325").Value = -1
$ws.Range("AA325").Value = 0.875
$ws.Range("AB325").Value = -1
$ws.Range("AC325").Value = 0.925

# Row 326
$ws.Range("B326").Value = 7302795
$ws.Range("F326").Value = "Unin Comercio"
$ws.Range("G326").Value = "Deportivo Garcilaso"
$ws.Range("I326").Value = 2
$ws.Range("J326").Value = "A"
$ws.Range("K326").Value = 2.25
$ws.Range("L326").Value = 3.3
$ws.Range("M326").Value = 2.7
$ws.Range("N326").Value = 1.75
$ws.Range("O326").Value = 3.6
$ws.Range("P326").Value = 4
$ws.Range("Q326").Value = -0.5
$ws.Range("R326").Value = 1.8
$ws.Range("S326").Value = 2
$ws.Range("T326").Value = 2.75
$ws.Range("U326").Value = 1.825
$ws.Range("V326").Value = 1.975
$ws.Range("W326").Value = -1
$ws.Range("Y326").Value = 3
$ws.Range("AA326").Value = 1
$ws.Range("AB326").Value = 0.4125
$ws.Range("AC326").Value = -0.5

# Row 334
$ws.Range("B334").Value = 7384626
$ws.Range("F334").Value = "Sporting Cristal"
$ws.Range("G334").Value = "Alianza Atletico"
$ws.Range("H334").Value = 3
$ws.Range("J334").Value = "H"
$ws.Range("K334").Value = 1.3
$ws.Range("L334").Value = 5
$ws.Range("M334").Value = 9
$ws.Range("N334").Value = 1.166
$ws.Range("O334").Value = 6.5
$ws.Range("P334").Value = 13
$ws.Range("Q334").Value = -2
$ws.Range("R334").Value = 1.85
$ws.Range("S334").Value = 1.95
$ws.Range("T334").Value = 3.25
$ws.Range("U334").Value = 2
$ws.Range("V334").Value = 1.8
$ws.Range("W334").Value = 0.1659999999999999
$ws.Range("X334").Value = -1
$ws.Range("Z334").Value = 0.8500000000000001
$ws.Range("AA334").Value = -1
$ws.Range("AB334").Value = -0.5
$ws.Range("AC334").Value = 0.4

# Row 336
$ws.Range("B336").Value = 7384625
$ws.Range("F336").Value = "AD Tarma"
$ws.Range("G336").Value = "Carlos Manucci"
$ws.Range("H336").Value = 0
$ws.Range("J336").Value = "D"
$ws.Range("K336").Value = 1.5
$ws.Range("L336").Value = 3.75
$ws.Range("M336").Value = 7
$ws.Range("N336").Value = 1.363
$ws.Range("O336").Value = 4.333
$ws.Range("P336").Value = 9.5
$ws.Range("Q336").Value = -1.25
$ws.Range("R336").Value = 1.875
$ws.Range("S336").Value = 1.925
$ws.Range("T336").Value = 2.5
$ws.Range("U336").Value = 1.8
$ws.Range("V336").Value = 2
$ws.Range("W336").Value = -1
$ws.Range("X336").Value = 3.333
$ws.Range("Z336").Value = -1
$ws.Range("AA336").Value = 0.925
$ws.Range("AB336").Value = -1
$ws.Range("AC336").Value = 1

# Row 369
$ws.Range("N369").Value = 1.444
$ws.Range("O369").Value = 4.333
$ws.Range("P369").Value = 6
$ws.Range("Q369").Value = -1.25
$ws.Range("R369").Value = 2.025
$ws.Range("S369").Value = 1.825
$ws.Range("T369").Value = 2.75
$ws.Range("U369").Value = 2.05
$ws.Range("V369").Value = 1.8

# Row 370
$ws.Range("N370").Value = 2.4
$ws.Range("O370").Value = 3.1
$ws.Range("P370").Value = 2.9
$ws.Range("Q370").Value = -0.25
$ws.Range("R370").Value = 2.125
$ws.Range("S370").Value = 1.75
$ws.Range("U370").Value = 2.025
$ws.Range("V370").Value = 1.825

# Row 371
$ws.Range("O371").Value = 3.6
$ws.Range("P371").Value = 6.5
$ws.Range("R371").Value = 2
$ws.Range("S371").Value = 1.85
$ws.Range("T371").Value = 2.5
$ws.Range("U371").Value = 2.05
$ws.Range("V371").Value = 1.8

# Row 372
$ws.Range("N372").Value = 4.75
$ws.Range("O372").Value = 3.6
$ws.Range("P372").Value = 1.7
$ws.Range("R372").Value = 1.925
$ws.Range("S372").Value = 1.925
$ws.Range("U372").Value = 1.85
$ws.Range("V372").Value = 2

# Row 373
$ws.Range("P373").Value = 7.5
$ws.Range("R373").Value = 2.05
$ws.Range("S373").Value = 1.8
$ws.Range("U373").Value = 2.025
$ws.Range("V373").Value = 1.825

# Row 374
$ws.Range("N374").Value = 1.8
$ws.Range("P374").Value = 4.2
$ws.Range("Q374").Value = -0.75
$ws.Range("R374").Value = 2.025
$ws.Range("S374").Value = 1.825
$ws.Range("U374").Value = 1.9
$ws.Range("V374").Value = 1.95

# Row 375
$ws.Range("N375").Value = 1.666
$ws.Range("P375").Value = 4.75
$ws.Range("R375").Value = 1.925
$ws.Range("S375").Value = 1.925
$ws.Range("U375").Value = 2
$ws.Range("V375").Value = 1.85
